$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization typos in existing shared strings
$ws.Range("C7").Value = "Entomology, Experimental Design, Material Science"
$ws.Range("C9").Value = "Experimental Design, Forensics"

# Add new event "Microbe Mission" as a new row in the events list (column F),
# copying the formatting/style of the previous row so it reuses the same style.
$ws.Range("F23").Copy($ws.Range("F24"))
$ws.Range("F24").Value = "Microbe Mission"
